$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column P: "mapping" header (row1) + "NA" values (rows 2-3) ---
# P1 reuses the same formatting already present on X1 (bold font + left/right
# thin border) so the engine dedups onto the existing style record instead of
# minting a throwaway one.
$ws.Range("X1").Copy()
$ws.Range("P1").PasteSpecial(-4122)
$ws.Range("P1").Value = "mapping"

# P2 gets a left/right thin border (same border used elsewhere on the sheet)
# without the bold font.
$ws.Range("P2").Value = "NA"
$leftBorder = $ws.Range("P2").Borders.Item(7)
$leftBorder.ColorIndex = 1
$leftBorder.LineStyle = 1
$rightBorder = $ws.Range("P2").Borders.Item(10)
$rightBorder.ColorIndex = 1
$rightBorder.LineStyle = 1

# P3 copies P2's formatting (same style reused) and gets the same value.
$ws.Range("P2").Copy()
$ws.Range("P3").PasteSpecial(-4122)
$ws.Range("P3").Value = "NA"

$excel.CutCopyMode = 0

# --- Update the remembered selection (cosmetic cursor position) ---
$ws.Range("P18").Select()
